$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 53.857143
$ws.Range("I11").Value = 53.857143
$ws.Range("K11").Value = 53.857143
$ws.Range("M11").Value = 86.14285699999999
$ws.Range("H33").Value = 200.18182
$ws.Range("I33").Value = 264.33334
$ws.Range("J33").Value = 155.76923
$ws.Range("K33").Value = 264.33334
$ws.Range("L33").Value = 155.76923
$ws.Range("M33").Value = -35.33334000000002
$ws.Range("N33").Value = -613.76923
$ws.Range("H46").Value = 909.5
$ws.Range("I46").Value = 886.3333
$ws.Range("K46").Value = 2658.9999
$ws.Range("M46").Value = -2539.9999
$ws.Range("H59").Value = 1005
$ws.Range("I59").Value = 900
$ws.Range("J59").Value = 1057.5
$ws.Range("K59").Value = 2700
$ws.Range("L59").Value = 3172.5
$ws.Range("M59").Value = -2143
$ws.Range("N59").Value = -4286.5
$ws.Range("H60").Value = 909.5
$ws.Range("I60").Value = 886.3333
$ws.Range("K60").Value = 2658.9999
$ws.Range("M60").Value = -2174.9999
$ws.Range("H64").Value = 10251.5
$ws.Range("I64").Value = 5003
$ws.Range("J64").Value = 12001
$ws.Range("K64").Value = 5003
$ws.Range("L64").Value = 12001
$ws.Range("M64").Value = -4755
$ws.Range("N64").Value = -12497
$ws.Range("H67").Value = 10251.5
$ws.Range("I67").Value = 5003
$ws.Range("J67").Value = 12001
$ws.Range("K67").Value = 5003
$ws.Range("L67").Value = 12001
$ws.Range("M67").Value = -4145
$ws.Range("N67").Value = -13717
$ws.Range("H86").Value = 52631576
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 52631576
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 52631576
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -52633822
$ws.Range("H89").Value = 52631576
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 52631576
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 263157880
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -263169112
$ws.Range("H106").Value = 3626.5454
$ws.Range("I106").Value = 5579.6
$ws.Range("K106").Value = 5579.6
$ws.Range("M106").Value = -4948.6
$ws.Range("H107").Value = 92651.27
$ws.Range("I107").Value = 112640.445
$ws.Range("K107").Value = 112640.445
$ws.Range("M107").Value = -110720.445
$ws.Range("H132").Value = 8255.799999999999
$ws.Range("I132").Value = 11152.272
$ws.Range("K132").Value = 33456.81600000001
$ws.Range("M132").Value = -30926.81600000001
$ws.Range("H138").Value = 4964.617
$ws.Range("I138").Value = 1554.8667
$ws.Range("J138").Value = 6562.9375
$ws.Range("K138").Value = 4664.6001
$ws.Range("L138").Value = 19688.8125
$ws.Range("M138").Value = 475.3999000000003
$ws.Range("N138").Value = -29968.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4132.0757
$ws.Range("I32").Value = 4050
$ws.Range("K32").Value = 4050
$ws.Range("M32").Value = -3763
$ws.Range("H61").Value = 3344.4546
$ws.Range("I61").Value = 1559.0869
$ws.Range("J61").Value = 7450.8
$ws.Range("K61").Value = 1559.0869
$ws.Range("L61").Value = 7450.8
$ws.Range("M61").Value = -1347.0869
$ws.Range("N61").Value = -7874.8
$ws.Range("H109").Value = 80000
$ws.Range("J109").Value = 80000
$ws.Range("L109").Value = 80000
$ws.Range("N109").Value = -82774
$ws.Range("H136").Value = 3344.4546
$ws.Range("I136").Value = 1559.0869
$ws.Range("J136").Value = 7450.8
$ws.Range("K136").Value = 4677.2607
$ws.Range("L136").Value = 22352.4
$ws.Range("M136").Value = -2127.2607
$ws.Range("N136").Value = -27452.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 856.6087
$ws.Range("I64").Value = 827
$ws.Range("J64").Value = 869.5625
$ws.Range("K64").Value = 827
$ws.Range("L64").Value = 869.5625
$ws.Range("M64").Value = -602
$ws.Range("N64").Value = -1319.5625
$ws.Range("H67").Value = 856.6087
$ws.Range("I67").Value = 827
$ws.Range("J67").Value = 869.5625
$ws.Range("K67").Value = 827
$ws.Range("L67").Value = 869.5625
$ws.Range("M67").Value = -47
$ws.Range("N67").Value = -2429.5625
$ws.Range("H107").Value = 986.25
$ws.Range("I107").Value = 986.25
$ws.Range("K107").Value = 986.25
$ws.Range("M107").Value = 933.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 500
$ws.Range("I8").Value = 500
$ws.Range("K8").Value = 500
$ws.Range("M8").Value = -360
$ws.Range("H62").Value = 3098
$ws.Range("I62").Value = 3098
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3098
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2474
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 3098
$ws.Range("I65").Value = 3098
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 15490
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -12370
$ws.Range("N65").ClearContents()
$ws.Range("H99").Value = 5727.091
$ws.Range("I99").Value = 4050
$ws.Range("K99").Value = 4050
$ws.Range("M99").Value = -2552
$ws.Range("H107").Value = 1856.4286
$ws.Range("I107").Value = 1423.4117
$ws.Range("J107").Value = 3696.75
$ws.Range("K107").Value = 1423.4117
$ws.Range("L107").Value = 3696.75
$ws.Range("M107").Value = 496.5882999999999
$ws.Range("N107").Value = -7536.75
$ws.Range("H126").Value = 5727.091
$ws.Range("I126").Value = 4050
$ws.Range("K126").Value = 12150
$ws.Range("M126").Value = -9680
$ws.Range("H132").Value = 4662.6875
$ws.Range("I132").Value = 2944.1428
$ws.Range("J132").Value = 5999.3335
$ws.Range("K132").Value = 8832.428400000001
$ws.Range("L132").Value = 17998.0005
$ws.Range("M132").Value = -6302.428400000001
$ws.Range("N132").Value = -23058.0005
$ws.Range("H134").Value = 4410.1787
$ws.Range("I134").Value = 3588.4443
$ws.Range("J134").Value = 5889.3
$ws.Range("K134").Value = 10765.3329
$ws.Range("L134").Value = 17667.9
$ws.Range("M134").Value = -8230.332900000001
$ws.Range("N134").Value = -22737.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2914559.2
$ws.Range("I4").Value = 3158012
$ws.Range("J4").Value = 845210
$ws.Range("K4").Value = 9474036
$ws.Range("L4").Value = 2535630
$ws.Range("M4").Value = -9473924
$ws.Range("N4").Value = -2535854
$ws.Range("H23").Value = 295.44446
$ws.Range("I23").Value = 227.4
$ws.Range("K23").Value = 682.2
$ws.Range("M23").Value = -447.2
$ws.Range("H75").Value = 703.5
$ws.Range("J75").Value = 695
$ws.Range("L75").Value = 2085
$ws.Range("N75").Value = -4081
$ws.Range("H78").Value = 703.5
$ws.Range("J78").Value = 695
$ws.Range("L78").Value = 6255
$ws.Range("N78").Value = -16239
$ws.Range("H117").Value = 1966.6666
$ws.Range("I117").Value = 1900
$ws.Range("J117").Value = 2000
$ws.Range("K117").Value = 5700
$ws.Range("L117").Value = 6000
$ws.Range("M117").Value = -2258
$ws.Range("N117").Value = -12884
$ws.Range("H121").Value = 1112033.2
$ws.Range("I121").Value = 399.83334
$ws.Range("J121").Value = 3335300
$ws.Range("K121").Value = 1199.50002
$ws.Range("L121").Value = 10005900
$ws.Range("M121").Value = 110.4999800000001
$ws.Range("N121").Value = -10008520

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1434736.8
$ws.Range("I80").Value = 2005499.8
$ws.Range("J80").Value = 915861.2
$ws.Range("K80").Value = 2005499.8
$ws.Range("L80").Value = 915861.2
$ws.Range("M80").Value = -2004501.8
$ws.Range("N80").Value = -917857.2
$ws.Range("H83").Value = 1434736.8
$ws.Range("I83").Value = 2005499.8
$ws.Range("J83").Value = 915861.2
$ws.Range("K83").Value = 10027499
$ws.Range("L83").Value = 4579306
$ws.Range("M83").Value = -10022507
$ws.Range("N83").Value = -4589290
$ws.Range("H132").Value = 1670599.9
$ws.Range("I132").Value = 2146487
$ws.Range("J132").Value = 4995
$ws.Range("K132").Value = 6439461
$ws.Range("L132").Value = 14985
$ws.Range("M132").Value = -6436931
$ws.Range("N132").Value = -20045

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5144.6
$ws.Range("I46").Value = 3908.3333
$ws.Range("K46").Value = 3908.3333
$ws.Range("M46").Value = -3720.3333
$ws.Range("H82").Value = 1896.4445
$ws.Range("I82").Value = 1734.6364
$ws.Range("K82").Value = 1734.6364
$ws.Range("M82").Value = -1373.6364
$ws.Range("H85").Value = 1896.4445
$ws.Range("I85").Value = 1734.6364
$ws.Range("K85").Value = 1734.6364
$ws.Range("M85").Value = -486.6364000000001
$ws.Range("H100").Value = 201920.4
$ws.Range("J100").Value = 2000
$ws.Range("L100").Value = 2000
$ws.Range("N100").Value = -3082
$ws.Range("H122").Value = 1170503.5
$ws.Range("I122").Value = 1112561.9
$ws.Range("J122").Value = 1228445.2
$ws.Range("K122").Value = 3337685.7
$ws.Range("L122").Value = 3685335.6
$ws.Range("M122").Value = -3335235.7
$ws.Range("N122").Value = -3690235.6
$ws.Range("H136").Value = 4841.1763
$ws.Range("I136").Value = 3279.7
$ws.Range("J136").Value = 7071.857
$ws.Range("K136").Value = 9839.099999999999
$ws.Range("L136").Value = 21215.571
$ws.Range("M136").Value = -7289.099999999999
$ws.Range("N136").Value = -26315.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1746.069
$ws.Range("I132").Value = 1523.5555
$ws.Range("J132").Value = 4750
$ws.Range("K132").Value = 4570.666499999999
$ws.Range("L132").Value = 14250
$ws.Range("M132").Value = -2040.666499999999
$ws.Range("N132").Value = -19310
$ws.Range("H136").Value = 315603.06
$ws.Range("I136").Value = 386850.12
$ws.Range("J136").Value = 6865.8335
$ws.Range("K136").Value = 1160550.36
$ws.Range("L136").Value = 20597.5005
$ws.Range("M136").Value = -1158000.36
$ws.Range("N136").Value = -25697.5005

Write-Output "Applied all cell updates"